$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 34.1829625
$ws.Range("H2").Value = 68.365925
$ws.Range("I2").Value = 0.6188383653293237
$ws.Range("J2").Value = 0.5689764244710266
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.36926145240625
$ws.Range("R2").Value = 1.477045809625
$ws.Range("S2").Value = 0.6188383653293237
$ws.Range("T2").Value = 0.5689764244710266

# Row 3
$ws.Range("I3").Value = 0.06738425137939692
$ws.Range("J3").Value = 0.09293230485581538
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 0.04020824811083333
$ws.Range("R3").Value = 0.241249488665
$ws.Range("S3").Value = 0.06738425137939692
$ws.Range("T3").Value = 0.09293230485581538

# Row 4
$ws.Range("G4").Value = 1.560491
$ws.Range("H4").Value = 4.681473
$ws.Range("I4").Value = 0.02825067311094296
$ws.Range("J4").Value = 0.03896162845449177
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.0108025
$ws.Range("N4").Value = 0.021605
$ws.Range("Q4").Value = 0.0168572040275
$ws.Range("R4").Value = 0.101143224165
$ws.Range("S4").Value = 0.02825067311094296
$ws.Range("T4").Value = 0.03896162845449177

# Row 5
$ws.Range("G5").Value = 11.3729585
$ws.Range("H5").Value = 22.745917
$ws.Range("I5").Value = 0.2058927176688748
$ws.Range("J5").Value = 0.1893032314851988
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0108025
$ws.Range("N5").Value = 0.021605
$ws.Range("Q5").Value = 0.12285638419625
$ws.Range("R5").Value = 0.491425536785
$ws.Range("S5").Value = 0.2058927176688748
$ws.Range("T5").Value = 0.1893032314851988

# Row 6
$ws.Range("G6").Value = 2.577819333333334
$ws.Range("H6").Value = 7.733458000000001
$ws.Range("I6").Value = 0.04666808800888241
$ws.Range("J6").Value = 0.0643618188686375
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.0108025
$ws.Range("N6").Value = 0.021605
$ws.Range("Q6").Value = 0.02784689334833334
$ws.Range("R6").Value = 0.16708136009
$ws.Range("S6").Value = 0.04666808800888241
$ws.Range("T6").Value = 0.0643618188686375

# Row 7
$ws.Range("G7").Value = 1.820947666666666
$ws.Range("H7").Value = 5.462842999999999
$ws.Range("I7").Value = 0.0329659045025792
$ws.Range("J7").Value = 0.04546459186482997
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.0108025
$ws.Range("N7").Value = 0.021605
$ws.Range("Q7").Value = 0.01967078716916666
$ws.Range("R7").Value = 0.118024723015
$ws.Range("S7").Value = 0.0329659045025792
$ws.Range("T7").Value = 0.04546459186482997
